$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width updates (OOXML stored width = COM ColumnWidth + 5/6) ---
$ws.Columns.Item(3).ColumnWidth = 33 - 5/6
$ws.Columns.Item(4).ColumnWidth = 84 - 5/6
$ws.Columns.Item(8).ColumnWidth = 30 - 5/6

# --- Row 2 ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1328547"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1328547"
$ws.Range("C2").Value = "SEO"
$ws.Range("D2").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("F2").Value = "10 applicants"
$ws.Range("H2").Value = "TAR - Company"

# --- Row 3 ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "1328541"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1328541"
$ws.Range("C3").Value = "Marketing"
$ws.Range("D3").Value = "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt"
$ws.Range("F3").Value = "18 applicants"
$ws.Range("H3").Value = "TAR - Company"

# --- Row 4 ---
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "1328155"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1328155"
$ws.Range("C4").Value = "Clinical Study Analyst Trainee"
$ws.Range("D4").Value = "Bruxelles, Belgio"
$ws.Range("F4").Value = "71 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "UCB"

# --- Row 5 ---
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "1327922"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1327922"
$ws.Range("C5").Value = "Digital Marketing Intern"
$ws.Range("D5").Value = "Nugegoda, Sri Lanka"
$ws.Range("F5").Value = "22 applicants"
$ws.Range("G5").Value = "3 - 6 Months"
$ws.Range("H5").Value = "Starbeans Ceylon (Pvt ) Ltd"

# --- Row 6 ---
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "1327889"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1327889"
$ws.Range("C6").Value = "Graphic Designer"
$ws.Range("D6").Value = "Birkat as SAB, Madinet Berkat as Sabee, Birket el Sab, Menofia Governorate, Egypt"
$ws.Range("F6").Value = "2 applicants"
$ws.Range("G6").Value = "9 - 12 Weeks"
$ws.Range("H6").Value = "Lines"
